$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates -----------------------------------------------------
# A3: data file name changed
$ws.Range("A3").Value = "ukb51139_subset.csv"

# B3: new "Dimensions" value added under the new data file
$ws.Range("B3").Value = "28012 x 1081"

# J3/K3/L3: new Vars Used / Overlap values added for row 3
$ws.Range("J3").Value = 91
$ws.Range("K3").Value = "94.0 & 81.9"
$ws.Range("L3").Value = "66.1 & 53.1"

# --- Column widths (match autosize/adjustment seen after edit) --------
# (columns 5 and 7 are unchanged by the edit, so they are left alone)
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 10.333333333333334
$ws.Columns.Item(8).ColumnWidth = 9.833333333333334
$ws.Columns.Item(9).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(10).ColumnWidth = 10.0
$ws.Columns.Item(11).ColumnWidth = 14.333333333333334
$ws.Columns.Item(12).ColumnWidth = 13.333333333333334

# --- Selection / active cell -------------------------------------------
$ws.Range("L4").Select()
